# Hjemme passive tweaks - lichtwark deleted values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header: replace subject IDs for columns B:E (1,2,3,4 -> 16,20,16,20)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 (CON): B2, D2, E2 deleted (no value for those subjects); C2 replaced
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = 5.0551662443276495
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()

# Row 3 (STR): replace with new values for the re-mapped subjects
$ws.Range("B3").Value = 5.3564084336275419
$ws.Range("C3").Value = 6.4263281534632748
$ws.Range("D3").Value = 8.0558901479131446
$ws.Range("E3").Value = 3.3256809647001098

# Update the selection to match the new authored state
$ws.Range("B1:E3").Select()
